$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'1120194100404 "
$ws.Range("B3").Value = "'1220194200661"
$ws.Range("B4").Value = "'0420194406627"
$ws.Range("B8").Value = "'1120170200928"
$ws.Range("B7").Value = "'0420172008461"
$ws.Range("B6").Value = "'1220170301419"
$ws.Range("B5").Value = "'1220170301418   "

$ws.Range("D8").Select()
